# Part_counts.xlsx update
#
# Source-of-truth context: Sheet1 tallies, per manufacturer, how many of
# each part type were dropped (via COUNTIFS against the cached rows of an
# external "parts" workbook link, [1]parts!...). The commit adds six more
# parts to the "Bright" manufacturer's drop log (2x pommel, 4x handle),
# which bumps:
#   - handle/Bright (row 4, col B) from 1 -> 5
#   - pommel/Bright (row 5, col B) from 1 -> 3
# The external link itself isn't reachable from this workbook (no real
# parts.csv / linked workbook is open in this session), so we push the new
# totals straight onto the dependent cells, same as Excel shows once a
# linked workbook has been refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated COUNTIFS results now that "Bright" has more handle/pommel drops.
$ws.Range("B4").Value = 5   # handle, Bright
$ws.Range("B5").Value = 3   # pommel, Bright

# Leftover UI state from the editing session (active cell moved to M16).
$ws.Range("M16").Select() | Out-Null
